$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Word-search board (rows 111-114) ---
$ws.Range("J111").Value = "o"
$ws.Range("K111").Value = "a"
$ws.Range("L111").Value = "b"
$ws.Range("M111").Value = "n"

$ws.Range("J112").Value = "o"
$ws.Range("K112").Value = "t"
$ws.Range("L112").Value = "a"
$ws.Range("M112").Value = "e"

$ws.Range("J113").Value = "a"
$ws.Range("K113").Value = "h"
$ws.Range("L113").Value = "k"
$ws.Range("M113").Value = "r"

$ws.Range("J114").Value = "a"
$ws.Range("K114").Value = "f"
$ws.Range("L114").Value = "l"
$ws.Range("M114").Value = "v"

# --- Trie / word list illustration (rows 118-121) ---
$ws.Range("J118").Value = "o"
$ws.Range("K118").Value = "a"
$ws.Range("L118").Value = "a"
$ws.Range("M118").Value = "n "

$ws.Range("J119").Value = "e"
$ws.Range("K119").Value = "t"
$ws.Range("L119").Value = "a"
$ws.Range("M119").Value = "e"

$ws.Range("J120").Value = "i"
$ws.Range("K120").Value = "h"
$ws.Range("L120").Value = "k"
$ws.Range("M120").Value = "r"

$ws.Range("J121").Value = "i"
$ws.Range("K121").Value = "f"
$ws.Range("L121").Value = "l"
$ws.Range("M121").Value = "v"

# --- Final note (row 123) ---
$ws.Range("I123").Value = "node"
$ws.Range("J123").Value = "o"

# --- Restore the active selection/view to match the edited area ---
$ws.Range("J119").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 115
$excel.ActiveWindow.ScrollColumn = 8
